$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark execution mode ("Manual") for every test case row in column E
$ws.Range("E2").Value = "Manual"
$ws.Range("E3").Value = "Manual"
$ws.Range("E4").Value = "Manual"
$ws.Range("E5").Value = "Manual"
$ws.Range("E6").Value = "Manual"
$ws.Range("E7").Value = "Manual"
$ws.Range("E8").Value = "Manual"
$ws.Range("E9").Value = "Manual"

# Row 5 wraps to a shorter height after the test case content update
$ws.Rows.Item(5).RowHeight = 195

# Move the active selection to E10 (just past the last data row)
$ws.Range("E10").Select()
